{"js": "// Updated manual docx to reflect README\n//\n// 1. Copyright line: \"2014-2018\" -> \"2014-2020\" and add \", Talitha Forcier\"\n//    between \"Eric Paniagua\" and \"Oliver Tam & Molly Hammell\" (with the\n//    auto \"_GoBack\" bookmark now sitting right before \"Oliver\").\n// 2. Requirements section: Python version note updated.\n// 3. The \"_GoBack\" bookmark (Word's \"last edit location\" marker) moves from\n//    its old spot (after \"...QIAseq stranded\") to the new edit location in\n//    the copyright line; every other bookmark keeps its name/position and\n//    just gets renumbered by Word as a side effect of that move.\n\nconst body = context.document.body;\n\n// --- Move the \"_GoBack\" bookmark ---------------------------------------\n// Word maintains a single \"_GoBack\" bookmark at the site of the most recent\n// edit. Remove the old one (its exact old position doesn't matter - Word\n// looks it up by name) before we insert the new one below.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- 1. Copyright paragraph ---------------------------------------------\nconst copyrightHits = body.search(\"Copyright (C)\", { matchCase: true });\ncopyrightHits.load(\"items\");\nawait context.sync();\n\nconst copyrightPara = copyrightHits.items[0].paragraphs.getFirst();\n\n// 2014-2018 -> 2014-2020\nconst yearHits = copyrightPara.search(\"2014-2018\", { matchCase: true });\nyearHits.load(\"items\");\nawait context.sync();\nyearHits.items[0].insertText(\"2014-2020\", \"Replace\");\nawait context.sync();\n\n// Insert the new \"Talitha Forcier, \" credit right before \"Oliver Tam & Molly\"\nconst oliverHits = copyrightPara.search(\"Oliver Tam & Molly\", { matchCase: true });\noliverHits.load(\"items\");\nawait context.sync();\n\nconst insertedRange = oliverHits.items[0].insertText(\"Talitha Forcier, \", \"Before\");\nawait context.sync();\n\n// Re-plant \"_GoBack\" as a zero-length bookmark right before \"Oliver\", i.e.\n// right after the text we just inserted.\nconst caret = insertedRange.getRange(\"After\");\ncaret.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- 2. Requirements > Python version note ------------------------------\nconst pyHits = body.search(\": 2.6.x or 2.7.x (not tested in Python 3.x)\", { matchCase: true });\npyHits.load(\"items\");\nawait context.sync();\npyHits.items[0].insertText(\": 2.6.x or 2.7.x or 3.x (tested on Python 2.7.11 and Python 3.7.7)\", \"Replace\");\nawait context.sync();\n", "ps1": "# Updated manual docx to reflect README\n#\n# 1. Copyright line: \"2014-2018\" -> \"2014-2020\" and add \", Talitha Forcier\"\n#    between \"Eric Paniagua\" and \"Oliver Tam & Molly Hammell\" (with Word's\n#    auto \"_GoBack\" bookmark re-planted right before \"Oliver\", at the new\n#    edit location).\n# 2. Requirements section: Python version note updated.\n\n$d = $word.ActiveDocument\n\n# --- 1. Copyright paragraph ---------------------------------------------\n$copyrightPara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"Copyright (C)*\") {\n        $copyrightPara = $p\n        break\n    }\n}\n\n# 2014-2018 -> 2014-2020\n$yearRange = $copyrightPara.Range.Duplicate\n$yearRange.Find.Execute(\"2014-2018\", $true, $false, $false, $false, $false, $true, 1, $false, \"2014-2020\", 2)\n\n# Insert the new \"Talitha Forcier, \" credit right before \"Oliver Tam & Molly\"\n$oliverRange = $copyrightPara.Range.Duplicate\n$oliverRange.Find.Execute(\"Oliver Tam & Molly\", $true)\n$oliverRange.InsertBefore(\"Talitha Forcier, \")\n\n# Re-find \"Oliver Tam & Molly\" (range shape settled after the insert) and\n# collapse to a zero-length caret right before it, then re-plant \"_GoBack\"\n# there - Word keeps a single \"_GoBack\" bookmark at the most recent edit\n# location, so adding one under that name moves the existing one.\n$oliverCaret = $copyrightPara.Range.Duplicate\n$oliverCaret.Find.Execute(\"Oliver Tam & Molly\", $true)\n$oliverCaret.Collapse(1)   # wdCollapseStart\n$d.Bookmarks.Add(\"_GoBack\", $oliverCaret)\n\n# --- 2. Requirements > Python version note ------------------------------\n$pyRange = $d.Content\n$pyRange.Find.Execute(\": 2.6.x or 2.7.x (not tested in Python 3.x)\", $true, $false, $false, $false, $false, $true, 1, $false, \": 2.6.x or 2.7.x or 3.x (tested on Python 2.7.11 and Python 3.7.7)\", 2)\n"}
